$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting for rule R10 (row 8) from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Select cell E8, matching the active selection recorded in the saved sheet view
$ws.Range("E8").Select()
